$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "CytPix" derived size-range values for column I (rows 2-19),
# plus the updated header in I1.
$values = @{
    1  = "CytPixSize"
    2  = "7 to 13"
    3  = "17 to >72"
    4  = "14 to 37"
    5  = "16 to >72"
    6  = "8 to 14"
    7  = "8 to 10"
    8  = "6 to 7"
    9  = "8 to 15"
    10 = "9 to 17"
    11 = "13 to >72"
    12 = "4 to 20"
    13 = "1 to 3<sup>8</sup>"
    14 = "6 to 9"
    15 = "5 to 10"
    16 = "7 to 31"
    17 = "12 to 15"
    18 = "6 to 8"
    19 = "10 to 17"
}

# Column I (rows 2-19) previously carried an accidental date number format
# (styles s="5"/s="6"). Re-style those cells to match the plain text style
# already used elsewhere in the table (column G, style index 2: Times New
# Roman 10pt, General format) before writing the new values, so Excel's
# "looks like a date" auto-detection doesn't reapply a date format.
for ($r = 2; $r -le 19; $r++) {
    $ws.Range("I$r").ClearFormats()
    $ws.Range("G$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
    $ws.Range("I$r").Value = $values[$r]
}

# Header cell I1 keeps its existing bold style; only its text changes.
$ws.Range("I1").Value = $values[1]

# Restore the selection to match the post-edit state.
$ws.Range("I13").Select()
